{"js": "// 1) Update \"Curso (semestre ideal): EB (5)\" -> \"Curso (semestre ideal): EM (10), EB (6)\"\nconst oldCurso = \"Curso (semestre ideal): EB (5)\";\nconst newCurso = \"Curso (semestre ideal): EM (10), EB (6)\";\n\nconst cursoHits = context.document.body.search(oldCurso, { matchCase: true });\ncursoHits.load(\"text\");\nawait context.sync();\n\nif (cursoHits.items.length > 0) {\n  cursoHits.items[0].insertText(newCurso, \"Replace\");\n  await context.sync();\n}\n\n// 2) Remove the \"Requisitos\" heading paragraph and the requirement bullet\n//    paragraph that follows it (\"LOT2059 -  Qu\u00edmica Org\u00e2nica Fundamental  (Requisito fraco)\").\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst toDelete = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text === \"Requisitos\" || text.indexOf(\"LOT2059\") !== -1) {\n    toDelete.push(paragraphs.items[i]);\n  }\n}\n\nfor (const p of toDelete) {\n  p.delete();\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Update \"Curso (semestre ideal): EB (5)\" -> \"Curso (semestre ideal): EM (10), EB (6)\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Execute(\"Curso (semestre ideal): EB (5)\", $false, $false, $false, $false, $false, $true, 1, $false, \"Curso (semestre ideal): EM (10), EB (6)\", 2) | Out-Null\n\n# 2) Remove the \"Requisitos\" heading paragraph and the requirement bullet\n#    paragraph that follows it (\"LOT2059 -  Qu\u00edmica Org\u00e2nica Fundamental  (Requisito fraco)\").\n$reqStart = $null\n$reqEnd = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $text = $p.Range.Text.Trim()\n    if ($text -eq \"Requisitos\") {\n        $reqStart = $p.Range.Start\n    }\n    if ($reqStart -ne $null -and $text.Contains(\"LOT2059\")) {\n        $reqEnd = $p.Range.End\n        break\n    }\n}\n\nif ($reqStart -ne $null -and $reqEnd -ne $null) {\n    $rng = $d.Range($reqStart, $reqEnd)\n    $rng.Delete()\n}\n"}
